$wb = $excel.ActiveWorkbook

$wsImported = $wb.Worksheets.Item("ImportedOntologies")
$wsConcepts = $wb.Worksheets.Item("Concepts")

# ---------------------------------------------------------------------------
# Concepts sheet (sheet3): new example rows 25 and 26 showing how to use a
# prefix in Manchester notation.
# Cells are written in the same order the shared strings were introduced so
# that the shared string table lines up with the source document.
# ---------------------------------------------------------------------------
$wsConcepts.Range("A25").Value = "SpecialMolecule"
$wsConcepts.Range("C25").Value = "Our own special molecules"
$wsConcepts.Range("D25").Value = "Used for our own special purpose"
$wsConcepts.Range("I25").Value = "hasPart some Atom"

# ---------------------------------------------------------------------------
# ImportedOntologies sheet (sheet2): new "prefix" column.
# ---------------------------------------------------------------------------
$wsImported.Range("B1").Value = "prefix"
$wsImported.Range("B3").Value = "emmo"

$wsConcepts.Range("A26").Value = "AnotherSpecialMolecule"
$wsConcepts.Range("G26").Value = "Molecule"
$wsConcepts.Range("G25").Value = "emmo-inferred-chemistry2:Molecule"
$wsConcepts.Range("I26").Value = "emmo-inferred-chemistry2:hasPart  some emmo-inferred-chemistry2:Atom"
$wsConcepts.Range("J25").Value = "Test giving prefix from emmo for subclass of and relations. Will be changed to emmo: when prefix is fixed"
$wsConcepts.Range("J26").Value = "Test giving prefix to relations. Will be changed to emmo: when prefix is fixed"

# ---------------------------------------------------------------------------
# ImportedOntologies sheet (sheet2): new "base_iri_root" column.
# ---------------------------------------------------------------------------
$wsImported.Range("C1").Value = "base_iri_root"
$wsImported.Range("C2").Value = "If base_iri_root is given, all imported ontologies whose base_iri starts with base_iri_root will be given the same prefix."
$wsImported.Range("B2").Value = "Local prefix for the imported ontology."

# Copy the header/description formatting onto the new columns.
$wsImported.Range("B1").Style = $wsImported.Range("A1").Style
$wsImported.Range("C1").Style = $wsImported.Range("A1").Style
$wsImported.Range("B2").Style = $wsImported.Range("A2").Style
$wsImported.Range("C2").Style = $wsImported.Range("A2").Style

# Widen the description row so the new long text is readable, and size the
# new columns.
$wsImported.Rows.Item(2).RowHeight = 75
$wsImported.Columns.Item(2).ColumnWidth = 17.28515625
$wsImported.Columns.Item(3).ColumnWidth = 28.85546875

# The imported-ontology example URL is now a real hyperlink.
$wsImported.Hyperlinks.Add($wsImported.Range("A3"), "https://raw.githubusercontent.com/emmo-repo/emmo-repo.github.io/master/versions/1.0.0-beta/emmo-inferred-chemistry2.ttl")

# An extra (empty) example row, styled like the hyperlink cell, plus a blank
# example value further down re-using the existing "space" shared string.
$wsImported.Range("C4").Style = "Hyperlink"
$wsImported.Range("D5").Value = " "

# ---------------------------------------------------------------------------
# Cosmetic formatting clean-up to match the final, re-saved layout.
# ---------------------------------------------------------------------------
$wsConcepts.Columns.Item(9).ColumnWidth = 68.6
$wsConcepts.Rows.Item(15).AutoFit()

# Restore the view state (active sheet/selection) seen in the saved file.
$wsConcepts.Activate()
$wsConcepts.Range("H30").Select()

$wsImported.Activate()
$wsImported.Range("A16").Select()
